# Commit: "Fruta / hortaliza, semanal"
#
# A new weekly price record was inserted as row 19 of the sheet (pushing the
# previously existing rows 19-77 down by one, to rows 20-78). The new record
# shares most attributes with its neighbours (same market, region, category,
# quality, volume, unit size, etc.) but carries a new date and new price
# figures for a "Perfection" variety sourced from "Provincia de Limarí".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 19; everything below (old rows 19-77)
# shifts down to rows 20-78, carrying its formatting and values with it.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44838
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112022
$ws.Cells.Item(19, 7).Value = "Arveja Verde"
$ws.Cells.Item(19, 8).Value = "Perfection"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 60
$ws.Cells.Item(19, 11).Value = 25000
$ws.Cells.Item(19, 12).Value = 26000
$ws.Cells.Item(19, 13).Value = 25500
$ws.Cells.Item(19, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 16).Value = 1020
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"
